# Lesson_Names.xlsx update
# - Insert a new row (row 36) with "Éducation aux médias" / "Médias"
#   right before "Éducation musicale", keeping the alphabetically sorted
#   table intact (everything below shifts down by one row).
# - Append a new row at the end of the table (row 82) with
#   "Sociologie" / "Socio".
# - Update the active view / selection to reflect where the user ended
#   up working (near the bottom of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new "Éducation aux médias" row above the current row 36
# ("Éducation musicale"); this pushes rows 36-80 down to 37-81.
$ws.Rows.Item(36).Insert()
$ws.Cells.Item(36, 2).Value2 = "Éducation aux médias"
$ws.Cells.Item(36, 3).Value2 = "Médias"

# Append the new "Sociologie" row after the (now shifted) last data
# row, 81.
$ws.Cells.Item(82, 2).Value2 = "Sociologie"
$ws.Cells.Item(82, 3).Value2 = "Socio"

# Column C carries the "Standard"-font style (style index 2) on every
# data row; copy it from the row above so the newly appended row
# matches the rest of the table instead of staying unformatted.
$ws.Cells.Item(81, 3).Copy()
$ws.Cells.Item(82, 3).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the final cursor/viewport position used while editing.
$win = $excel.ActiveWindow
try { $win.ScrollRow = 76 } catch {}
try { $win.ScrollColumn = 1 } catch {}
$ws.Range("B84").Select()
